$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 64, shifting existing rows 64..135 down to 65..136.
$ws.Rows.Item(64).Insert()

# Populate the newly inserted row 64 with the new price record.
$ws.Range("A64").Value = 2
$ws.Range("B64").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C64").Value = "Coquimbo"
$ws.Range("D64").Value = 44874
$ws.Range("E64").Value = 4
$ws.Range("F64").Value = 100112024
$ws.Range("G64").Value = "Choclo"
$ws.Range("H64").Value = "Dulce o Americano"
$ws.Range("I64").Value = "Primera"
$ws.Range("J64").Value = 700
$ws.Range("K64").Value = 25000
$ws.Range("L64").Value = 27000
$ws.Range("M64").Value = 26000
$ws.Range("N64").Value = "$/malla 70 unidades"
$ws.Range("O64").Value = "Provincia de Limarí"
$ws.Range("P64").Value = 371
$ws.Range("Q64").Value = 70
$ws.Range("R64").Value = "Hortaliza"
